$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 684, shifting existing rows 684-704 down to 685-705
$ws.Rows.Item(684).Insert()

# Populate the newly inserted row 684 with the new weekly record
$ws.Cells.Item(684, 1).Value = 8
$ws.Cells.Item(684, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(684, 3).Value = "Coquimbo"
$ws.Cells.Item(684, 4).Value = 45239
$ws.Cells.Item(684, 5).Value = 4
$ws.Cells.Item(684, 6).Value = 100112043
$ws.Cells.Item(684, 7).Value = "Pepino dulce"
$ws.Cells.Item(684, 8).Value = "Sin especificar"
$ws.Cells.Item(684, 9).Value = "Segunda"
$ws.Cells.Item(684, 10).Value = 200
$ws.Cells.Item(684, 11).Value = 19000
$ws.Cells.Item(684, 12).Value = 20000
$ws.Cells.Item(684, 13).Value = 19500
$ws.Cells.Item(684, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(684, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(684, 16).Value = 1083
$ws.Cells.Item(684, 17).Value = 18
$ws.Cells.Item(684, 18).Value = "Hortaliza"
